$d = $word.ActiveDocument

# wdHeaderFooterIndex constants
$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2
$wdHeaderFooterEvenPages = 3

foreach ($sec in $d.Sections) {
    foreach ($idx in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage, $wdHeaderFooterEvenPages)) {
        $hdr = $sec.Headers.Item($idx)
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("Freude, Bild, Weg, Psalm23", $true, $false, $false, $false, $false,
                                     $true, 1, $false, "Freude, Bild, Weg", 2) | Out-Null
        }
    }
}
